# Trade #2 closed at 2026-02-17 12:26:27 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#   - Summary sheet: refresh aggregate metrics after the new trade closed
#   - Strategy Status sheet: refresh the MarketMaking strategy row
#   - All Trades / MarketMaking sheets: append the newly closed trade (#2)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.99   # Current Capital
$summary.Range("B4").Value = -0.01     # Total P&L $
$summary.Range("B5").Value = -0.1      # Total P&L %
$summary.Range("B6").Value = 2         # Total Trades
$summary.Range("B8").Value = 1         # Losing Trades
$summary.Range("B9").Value = 50        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.98999999999999  # Capital
$status.Range("D4").Value = 2                  # Trades
$status.Range("E4").Value = -0.01              # P&L $
$status.Range("F4").Value = -0.01              # P&L %
$status.Range("G4").Value = 50                 # Win Rate %

# ---------------------------------------------------------------------
# Helper: append the newly closed Trade #2 row to a trades sheet
# ---------------------------------------------------------------------
function Add-TradeTwoRow($ws) {
    $ws.Range("A3").Value = 2

    # "2026-02-17" looks like a date literal, so Excel would otherwise
    # auto-convert it to a date serial number on assignment. Force the
    # cell to text first so it is stored as the literal string, matching
    # the existing Date column cells.
    $ws.Range("B3").NumberFormat = "@"
    $ws.Range("B3").Value = "2026-02-17"

    $ws.Range("C3").Value = "12:26:20"
    $ws.Range("D3").Value = "MarketMaking"
    $ws.Range("E3").Value = "UP"
    $ws.Range("F3").Value = 0.22
    $ws.Range("G3").Value = 0.19
    $ws.Range("H3").Value = "CLOSED"
    $ws.Range("I3").Value = -13.6364
    $ws.Range("J3").Value = -0.03
    $ws.Range("K3").Value = 99.98999999999999
    $ws.Range("L3").Value = 0
    $ws.Range("M3").Value = 0
    $ws.Range("N3").Value = 0.6
    $ws.Range("O3").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P3").Value = "early_exit"
    $ws.Range("Q3").Value = 0.13
}

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeTwoRow $allTrades

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeTwoRow $marketMaking
